$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Update the first three rows in place ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Insert 10 new rows right after row 3, each holding one of the
#     values that used to be packed (tab separated) into the later
#     summary rows. `Rows.Add(beforeRow)` inserts a new blank row
#     immediately *before* beforeRow, so we repeatedly target the next
#     open slot (row 4) and advance. ---
$newValues = @("99", "0.00003", "0.00005", "0.00004", "0.00000", "0.00004", "0.00004", "0.00004", "0.00372", "100.0")

$insertPos = 4
foreach ($val in $newValues) {
    $refRow = $t.Rows.Item($insertPos)
    $newRow = $t.Rows.Add($refRow)
    $t.Rows.Item($insertPos).Cells.Item(1).Range.Text = $val
    $insertPos = $insertPos + 1
}

# --- The three former multi-run (tab separated) rows now sit 10 rows
#     further down the table; collapse each back to a single value. ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "221"
